# Updated TPM values: NATMI LR-pair output table now includes the
# "ECs" target-cluster rows (previously only "FAPs" was present) and
# refreshed numeric weights throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ccl5", "Ccr4", "ECs", 2, 0.6666666666666666, 0.3410236666666666, 1.023071, 0.01850325494520333, 0.01850325494520333, 1, 0.3333333333333333, 0.04581866666666667, 0.137456, 0.4376255591461182, 0.4376255591461182, 0.01562524970844444, 0.140627247376, 0.008097497291417783, 0.008097497291417781),
    @("ECs", "Ccl5", "Ccr4", "FAPs", 2, 0.6666666666666666, 0.3410236666666666, 1.023071, 0.01850325494520333, 0.01850325494520333, 1, 0.3333333333333333, 0.05887966666666666, 0.176639, 0.5623744408538818, 0.5623744408538818, 0.02007935981877777, 0.180714238369, 0.01040575765378554, 0.01040575765378554),
    @("FAPs", "Ccl5", "Ccr4", "ECs", 3, 1, 3.842846333333334, 11.528539, 0.2085050756621187, 0.2085050756621187, 1, 0.3333333333333333, 0.04581866666666667, 0.137456, 0.4376255591461182, 0.4376255591461182, 0.1760740951982222, 1.584666856784, 0.09124715032143839, 0.09124715032143836),
    @("FAPs", "Ccl5", "Ccr4", "FAPs", 3, 1, 3.842846333333334, 11.528539, 0.2085050756621187, 0.2085050756621187, 1, 0.3333333333333333, 0.05887966666666666, 0.176639, 0.5623744408538818, 0.5623744408538818, 0.2262655111578889, 2.036389600421, 0.1172579253406803, 0.1172579253406803),
    @("MuSCs", "Ccl5", "Ccr4", "ECs", 2, 0.6666666666666666, 0.2092423333333333, 0.627727, 0.0113530661283407, 0.0113530661283407, 1, 0.3333333333333333, 0.04581866666666667, 0.137456, 0.4376255591461182, 0.4376255591461182, 0.009587204723555556, 0.086284842512, 0.004968391912437955, 0.004968391912437954),
    @("MuSCs", "Ccl5", "Ccr4", "FAPs", 2, 0.6666666666666666, 0.2092423333333333, 0.627727, 0.0113530661283407, 0.0113530661283407, 1, 0.3333333333333333, 0.05887966666666666, 0.176639, 0.5623744408538818, 0.5623744408538818, 0.01232011883922222, 0.110881069553, 0.006384674215902746, 0.006384674215902746),
    @("Resolving-Mac", "Ccl5", "Ccr4", "ECs", 3, 1, 14.03735666666667, 42.11207, 0.7616386032643372, 0.7616386032643372, 1, 0.3333333333333333, 0.04581866666666667, 0.137456, 0.4376255591461182, 0.4376255591461182, 0.6431729659911112, 5.78855669392, 0.3333125196208241, 0.333312519620824),
    @("Resolving-Mac", "Ccl5", "Ccr4", "FAPs", 3, 1, 14.03735666666667, 42.11207, 0.7616386032643372, 0.7616386032643372, 1, 0.3333333333333333, 0.05887966666666666, 0.176639, 0.5623744408538818, 0.5623744408538818, 0.8265148814144444, 7.43863393273, 0.4283260836435131, 0.4283260836435131)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
